$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The two "summary" rows (product "Баланс" and the highlighted
# section-separator row "Дивный вечер") traded places in the refreshed
# export, and each keeps its own distinct formatting (the separator
# row uses a highlighted fill / no border-right, the normal data row
# uses the regular bordered style with an integer number format).
# Swap the two rows' contents+formatting via a scratch range so the
# per-row formatting travels with the data, then clear the scratch
# range completely so it does not enlarge the sheet's used range.
# ------------------------------------------------------------------
$ws.Range("A56:B56").Copy($ws.Range("A200:B200"))
$ws.Range("A56:B56").ClearContents()
$ws.Range("A57:B57").Copy($ws.Range("A56:B56"))
$ws.Range("A57:B57").ClearContents()
$ws.Range("A200:B200").Copy($ws.Range("A57:B57"))
$ws.Range("A200:B200").Clear()

# Row 56 ("Баланс") now holds the old row 57 quantity (414); update it
# to the refreshed quantity. Row 57 ("Дивный вечер") keeps no quantity.
$ws.Range("B56").Value = 90


# --- Update product rows 2-112 with refreshed stock data ---

$ws.Range("A2").Value = "ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы"
$ws.Range("B2").Value = 81632
$ws.Range("A3").Value = "Чага (березовый гриб) 50г"
$ws.Range("B3").Value = 17570
$ws.Range("A4").Value = "Спорыш трава 50г"
$ws.Range("B4").Value = 11335
$ws.Range("A5").Value = "Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г"
$ws.Range("B5").Value = 3181
$ws.Range("A6").Value = "Солодка корни 50г"
$ws.Range("B6").Value = 31238
$ws.Range("A7").Value = "Чистотел трава 50г"
$ws.Range("B7").Value = 13635
$ws.Range("A8").Value = "Сенна листья 50г"
$ws.Range("B8").Value = 17849
$ws.Range("A9").Value = "Мать-и-мачеха листья 35г"
$ws.Range("B9").Value = 22948
$ws.Range("A10").Value = "Шиповник плоды низковитаминные 50г"
$ws.Range("B10").Value = 25485
$ws.Range("A11").Value = "Пижма цветки 75г"
$ws.Range("B11").Value = 16163
$ws.Range("A12").Value = "Полынь горькая трава 50г"
$ws.Range("B12").Value = 36512
$ws.Range("A13").Value = "Алтей корни 75г"
$ws.Range("B13").Value = 5204
$ws.Range("A14").Value = "Липа цветки 35г"
$ws.Range("B14").Value = 16792
$ws.Range("A15").Value = "Кукуруза столбики с рыльцами 40г"
$ws.Range("B15").Value = 25370
$ws.Range("A16").Value = "Дуба кора 75г"
$ws.Range("B16").Value = 64575
$ws.Range("A17").Value = "Сб. Грудной №4 50г"
$ws.Range("B17").Value = 38216
$ws.Range("A18").Value = "Мята перечная листья 50г"
$ws.Range("B18").Value = 24537
$ws.Range("A19").Value = "Брусника листья 50г"
$ws.Range("B19").Value = 16489
$ws.Range("A20").Value = "Эвкалипт прутовидный листья 75г"
$ws.Range("B20").Value = 27907
$ws.Range("A21").Value = "Багульник болотный побеги 50г"
$ws.Range("B21").Value = 15489
$ws.Range("A22").Value = "Ноготки цветки 50г"
$ws.Range("B22").Value = 27215
$ws.Range("A23").Value = "Ромашка цветки вн 50г"
$ws.Range("B23").Value = 106743
$ws.Range("A24").Value = "Береза почки 50г"
$ws.Range("B24").Value = 21259
$ws.Range("A25").Value = "Укроп пахучий плоды 50г"
$ws.Range("B25").Value = 73724
$ws.Range("A26").Value = "Чабрец трава 50г"
$ws.Range("B26").Value = 24165
$ws.Range("A27").Value = "Девясил корневища и корни 50г"
$ws.Range("B27").Value = 21163
$ws.Range("A28").Value = "Эрва шерстистая трава 30г"
$ws.Range("B28").Value = 16336
$ws.Range("A29").Value = "Пустырник трава 50г"
$ws.Range("B29").Value = 13916
$ws.Range("A30").Value = "Валериана корневища с корнями 50г"
$ws.Range("B30").Value = 24022
$ws.Range("A31").Value = "Боярышник плоды 75г"
$ws.Range("B31").Value = 26268
$ws.Range("A32").Value = "Сб. Фитонефрол (Урологический сбор) 50г"
$ws.Range("B32").Value = 10447
$ws.Range("A33").Value = "Подорожник большой листья 50г"
$ws.Range("B33").Value = 11424
$ws.Range("A34").Value = "Шалфей листья 50г"
$ws.Range("B34").Value = 44640
$ws.Range("A35").Value = "Бессмертник песчаный цветки 30г"
$ws.Range("B35").Value = 34499
$ws.Range("A36").Value = "Ламинарии слоевища (морская капуста) 100г"
$ws.Range("B36").Value = 21094
$ws.Range("A37").Value = "Аир корневища 75г"
$ws.Range("B37").Value = 10179
$ws.Range("A38").Value = "Лен семена 100г"
$ws.Range("B38").Value = 77829
$ws.Range("A39").Value = "Рябина плоды 50г"
$ws.Range("B39").Value = 2674
$ws.Range("A40").Value = "Крушина кора 50г"
$ws.Range("B40").Value = 15042
$ws.Range("A41").Value = "Череда трава 50г"
$ws.Range("B41").Value = 17517
$ws.Range("A42").Value = "Зверобой трава 50г"
$ws.Range("B42").Value = 55109
$ws.Range("A43").Value = "Можжевельник плоды 50г"
$ws.Range("B43").Value = 21278
$ws.Range("A44").Value = "Толокнянка листья 50г"
$ws.Range("B44").Value = 12529
$ws.Range("A45").Value = "Тысячелистник трава 50г"
$ws.Range("B45").Value = 24777
$ws.Range("A46").Value = "Сб. Фитопектол №1 (Грудной сбор №1) 35г"
$ws.Range("B46").Value = 9379
$ws.Range("A47").Value = "Крапива листья 50г"
$ws.Range("B47").Value = 26045
$ws.Range("A48").Value = "Сб. Фитопектол №2 (Грудной сбор №2) 35г"
$ws.Range("B48").Value = 12902
$ws.Range("A49").Value = "Фп Фиточай `"Лактафитол`" (БАД) 20х1,5 г"
$ws.Range("B49").Value = 13549
$ws.Range("A50").Value = "Фп Детский травяной чай `"ФармаЦветик® для иммунитета`" 20х1,5 г"
$ws.Range("B50").Value = 2738
$ws.Range("A51").Value = "Фп Детский травяной чай `"ФармаЦветик®  при простуде`" 20х1,5 г"
$ws.Range("B51").Value = 4481
$ws.Range("A52").Value = "Фп Детский травяной чай `"ФармаЦветик® для животика`" 20х1,5 г"
$ws.Range("B52").Value = 4700
$ws.Range("A53").Value = "Фп Детский травяной чай `"ФармаЦветик® для спокойного сна`" 20х1,5 г"
$ws.Range("B53").Value = 7638
$ws.Range("A54").Value = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем`"(БАД) 20*1,5г"
$ws.Range("B54").Value = 8260
$ws.Range("A55").Value = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем`" (БАД) 20*1,5г"
$ws.Range("B55").Value = 9750
$ws.Range("A58").Value = "Фп `"Щедрость природы® Фиточай для иммунитета`" 20х2,0 г"
$ws.Range("B58").Value = 1170
$ws.Range("A59").Value = "Фп `"Щедрость природы® Фиточай кардиологический`" 20х2,0 г"
$ws.Range("B59").Value = 1512
$ws.Range("A60").Value = "Фп `"Щедрость природы® Фиточай при простуде`" 20х2,0 г"
$ws.Range("B60").Value = 1134
$ws.Range("A61").Value = "Фп `"Щедрость природы® Фиточай успокоительный`"20х2,0 г"
$ws.Range("B61").Value = 2592
$ws.Range("A62").Value = "Фп Шалфей листья 20х1,5г"
$ws.Range("B62").Value = 113244
$ws.Range("A63").Value = "Фп Сб. Арфазетин-Э 20x2,0г"
$ws.Range("B63").Value = 28965
$ws.Range("A64").Value = "Фп `"Щедрость природы® Фиточай диабетический`" 20х2,0 г"
$ws.Range("B64").Value = 1098
$ws.Range("A65").Value = "Фп Сб. Грудной №4 20x2,0г"
$ws.Range("B65").Value = 566985
$ws.Range("A66").Value = "Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г"
$ws.Range("B66").Value = 62519
$ws.Range("A67").Value = "Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г"
$ws.Range("B67").Value = 45657
$ws.Range("A68").Value = "Фп Мята перечная листья 20x1,5г"
$ws.Range("B68").Value = 53453
$ws.Range("A69").Value = "Фп Подорожник листья 20x1,5г"
$ws.Range("B69").Value = 26425
$ws.Range("A70").Value = "Фп Сб. Бруснивер 20x2,0г"
$ws.Range("B70").Value = 177231
$ws.Range("A71").Value = "Фп Крапива листья 20x1,5г"
$ws.Range("B71").Value = 56215
$ws.Range("A72").Value = "Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г"
$ws.Range("B72").Value = 21304
$ws.Range("A73").Value = "Фп Липа цветки 20x1,5г"
$ws.Range("B73").Value = 65698
$ws.Range("A74").Value = "Фп Сб. Желудочный №3 20x2,0г"
$ws.Range("B74").Value = 23745
$ws.Range("A75").Value = "Фп Толокнянка листья 20x1,5г"
$ws.Range("B75").Value = 32922
$ws.Range("A76").Value = "Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г"
$ws.Range("B76").Value = 173371
$ws.Range("A77").Value = "Фп Аир корневища 20x1,5г"
$ws.Range("B77").Value = 3820
$ws.Range("A78").Value = "Фп Чабрец трава 20x1,5 г"
$ws.Range("B78").Value = 76029
$ws.Range("A79").Value = "Фп Ромашка цветки 20x1,5г"
$ws.Range("B79").Value = 1427339
$ws.Range("A80").Value = "Фп Мелисса лекарственная трава 20x1,5г"
$ws.Range("B80").Value = 40843
$ws.Range("A81").Value = "Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г"
$ws.Range("B81").Value = 83511
$ws.Range("A82").Value = "Фп Сб. Элекасол 20x2,0г"
$ws.Range("B82").Value = 44384
$ws.Range("A83").Value = "Фп Череда трава 20х1,5г"
$ws.Range("B83").Value = 50130
$ws.Range("A84").Value = "Фп Боярышник плоды 20х3,0г"
$ws.Range("B84").Value = 19196
$ws.Range("A85").Value = "Фп Сенна листья 20x1,5г"
$ws.Range("B85").Value = 75300
$ws.Range("A86").Value = "Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г"
$ws.Range("B86").Value = 95397
$ws.Range("A87").Value = "Фп Шиповник плоды 20х2,0г"
$ws.Range("B87").Value = 60887
$ws.Range("A88").Value = "Фп Фиточай `"Тибетский`" (БАД) 20х2,0  г"
$ws.Range("B88").Value = 9270
$ws.Range("A89").Value = "Фп Фиточай `"Опалиховский`" (БАД) 20х2,0 г"
$ws.Range("B89").Value = 5814
$ws.Range("A90").Value = "Фп Пижма цветки 20х1,5г"
$ws.Range("B90").Value = 6888
$ws.Range("A91").Value = "Фп `"Щедрость природы® Фиточай очищающий`" 20х2,0 г"
$ws.Range("B91").Value = 1854
$ws.Range("A92").Value = "Фп Зверобой трава 20x1,5г"
$ws.Range("B92").Value = 65689
$ws.Range("A93").Value = "Фп Брусника листья 20х1,5г"
$ws.Range("B93").Value = 94139
$ws.Range("A94").Value = "Фп Пустырник трава 20x1,5г"
$ws.Range("B94").Value = 53367
$ws.Range("A95").Value = "Фп Чистотел трава 20х1,5г"
$ws.Range("B95").Value = 39516
$ws.Range("A96").Value = "Фп `"Щедрость природы® Фиточай для пищеварения`" 20х2,0 г"
$ws.Range("B96").Value = 1890
$ws.Range("A97").Value = "Фп Душица трава 20x1,5г"
$ws.Range("B97").Value = 39222
$ws.Range("A98").Value = "Фп Пастушья сумка трава 20х1,5г"
$ws.Range("B98").Value = 7432
$ws.Range("A99").Value = "Фп Хвощ полевой трава 20х1,5г"
$ws.Range("B99").Value = 41920
$ws.Range("A100").Value = "Фп Береза листья 20x1,5г"
$ws.Range("B100").Value = 7240
$ws.Range("A101").Value = "Фп Золототысячник трава 20х1,5г"
$ws.Range("B101").Value = 7099
$ws.Range("A102").Value = "Фп Фиалка трехцветная трава 20x1,5г"
$ws.Range("B102").Value = 6664
$ws.Range("A103").Value = "Фп Ольха соплодия 20х1,5г"
$ws.Range("B103").Value = 6577
$ws.Range("A104").Value = "Фп Ноготки цветки 20x1,5г"
$ws.Range("B104").Value = 21720
$ws.Range("A105").Value = "Фп Кровохлебка корневища и корни 20x1,5г"
$ws.Range("B105").Value = 4390
$ws.Range("A106").Value = "Фп Почечный чай листья 20x1,5г"
$ws.Range("B106").Value = 54222
$ws.Range("A107").Value = "Фп Валериана корневища с корнями 20x1,5г"
$ws.Range("B107").Value = 19455
$ws.Range("A108").Value = "Фп Девясил корневища и корни 20х1,5г"
$ws.Range("B108").Value = 15882
$ws.Range("A109").Value = "Фп Лапчатка корневища 20x2,5г"
$ws.Range("B109").Value = 4240
$ws.Range("A110").Value = "Фп Тысячелистник трава 20x1,5г"
$ws.Range("B110").Value = 23758
$ws.Range("A111").Value = "Фп Крушина кора 20x1,5г"
$ws.Range("B111").Value = 13635
$ws.Range("A112").Value = "Фп Бадан корневища 20x1,5г"
$ws.Range("B112").Value = 2905
# --- Update the view: active cell / selection moved to A96 ---
$ws.Range("A96").Select()
